$d = $word.ActiveDocument

# --- Paragraph 1 formatting: add a (line-less) paragraph border with
#     5-twip spacing on every side, and widen the left indent from
#     120 (6pt) to 225 (11.25pt) twips. ---
$p1 = $d.Paragraphs(1)
$pf1 = $p1.Range.ParagraphFormat
$pf1.LeftIndent = 11.25

$borders1 = $pf1.Borders
$borders1.DistanceFromTop = 5
$borders1.DistanceFromLeft = 5
$borders1.DistanceFromBottom = 5
$borders1.DistanceFromRight = 5

# --- Replace the placeholder id text and drop the trailing space run
#     that followed it, leaving a single run with the new id text. ---
$range = $d.Content
$range.Find.Execute("**ID__AFFARS_5336_topic_2__ID** ", $false, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_SUBPART_5336_2__ID**", 2)
